$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): F1, G1, H1 get new labels; old I1:O1 removed
$ws.Range("F1").Value = "xgboost-y_drop"
$ws.Range("G1").Value = "lasso-y_drop"
$ws.Range("H1").Value = "num_datasets"

# Row 2 - reordered D2/E2 content, new F2/G2 arrays, H2 = num_datasets (14)
$ws.Range("D2").Value = "[[-0.02142273]`n [ 0.03085009]`n [-0.01658338]`n [-0.08774961]`n [ 0.07629745]`n [ 0.01285947]`n [-0.1510179 ]`n [ 0.12167721]`n [ 0.02495436]`n [-0.12302472]`n [-0.00984156]`n [ 0.14192943]`n [ 0.17040717]`n [-0.0973577 ]`n [ 0.09454775]`n [-0.13742305]`n [-0.10978463]`n [-0.0120889 ]]"
$ws.Range("E2").Value = "[[-0.13329466]`n [-0.0755709 ]`n [ 0.07038691]`n [-0.00164962]`n [ 0.03918755]`n [-0.02697289]`n [-0.1761752 ]`n [ 0.2074219 ]`n [-0.00913573]`n [-0.05056708]`n [-0.11520114]`n [ 0.21395498]`n [ 0.23900712]`n [-0.01852215]`n [ 0.17288869]`n [-0.17013782]`n [-0.13586598]`n [-0.12252524]]"
$ws.Range("F2").Value = "[-0.01534459  0.2239905   0.02075964  0.03055483  0.2244603  -0.00801081`n -0.0942831   0.0306251  -0.00213945 -0.00791237 -0.00733042  0.06793422`n -0.00980453 -0.00111156 -0.0084103  -0.18789811 -0.02736913  0.0831384 ]"
$ws.Range("G2").Value = "[-0.02183669  0.02984185 -0.00529348 -0.0746076   0.06522541  0.00451568`n -0.14802472  0.11230291  0.01416389 -0.10475166 -0.0128725   0.12858878`n  0.15872908 -0.08089117  0.08777961 -0.13037342 -0.10517074 -0.01009647]"
$ws.Range("H2").Value = 14

# Row 3 - reordered D3/E3 content, new F3/G3 arrays, H3 = num_datasets (14)
$ws.Range("D3").Value = "[[ 0.07821837]`n [-0.15734918]`n [ 0.02657613]`n [-0.00395006]`n [-0.03550372]`n [ 0.04690658]`n [ 0.04993961]`n [ 0.15074391]`n [ 0.01042418]`n [ 0.08316006]`n [-0.03055012]`n [-0.02010913]`n [ 0.06437487]`n [-0.17307925]`n [ 0.02442902]`n [-0.14350568]`n [-0.01515553]`n [-0.04834132]]"
$ws.Range("E3").Value = "[[ 0.11549848]`n [-0.16755368]`n [ 0.01758147]`n [-0.01034481]`n [-0.03462524]`n [ 0.05721805]`n [ 0.04972808]`n [ 0.14122358]`n [ 0.00536362]`n [ 0.07311682]`n [-0.04599812]`n [-0.0038856 ]`n [ 0.10697546]`n [-0.18829422]`n [ 0.02196575]`n [-0.15903066]`n [-0.01525848]`n [-0.05645174]]"
$ws.Range("F3").Value = "[ 0.0755531  -0.07964763 -0.03116396  0.01566319 -0.0019985   0.14138444`n  0.07730947  0.08924988 -0.0171617   0.07567372 -0.03096188 -0.0019985`n  0.05287367 -0.04099603  0.13360655 -0.12689601 -0.03096188  0.01770351]"
$ws.Range("G3").Value = "[ 9.27115254e-02 -1.65547540e-01  1.62093408e-02 -1.53449298e-02`n -3.52491453e-02  5.71114641e-02  5.04963199e-02  1.56787913e-01`n -5.70181680e-05  8.77986307e-02 -3.00183872e-02 -1.09736062e-02`n  7.87387679e-02 -1.82008811e-01  2.45032651e-02 -1.51574783e-01`n -5.74284815e-03 -6.06114024e-02]"
$ws.Range("H3").Value = 14

# Remove old columns I:O (this also shrinks the sheet dimension from O3 to H3)
$ws.Range("I1:O3").Clear()

# Restore default (non-custom) row heights after setting multi-line text values
$ws.Rows.Item(2).EntireRow.AutoFit()
$ws.Rows.Item(3).EntireRow.AutoFit()
